$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "admin"
$ws.Range("B3").Value = "8c6976e5b5410415bde908bd4dee15dfb167a9c873fc4bb8a81f6f2ab448a918"
$ws.Range("C3").Value = "admin@admin"
$ws.Range("D3").Value = "staff"
